$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8346506357192993
$ws.Range("B1").Value = 0.8278881311416626
$ws.Range("C1").Value = 1.354677319526672
$ws.Range("D1").Value = 2.390108585357666
$ws.Range("E1").Value = 1.771325707435608
